# Generate Report for Handback
# Updates the localization-status workbook to reflect that both language
# targets (zh-cn, de-de) have been handed back and are in sync with en-US,
# and records the resulting handback target/handoff-file/datetime info.

$wb = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (zh-cn = E, de-de = F) ---
$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText

# Overview status columns now hold the longer string -> widen columns
$ovw.Columns.Item(5).ColumnWidth = 29.9777047293527
$ovw.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/342d41c5bbd5653daaa7680fc7bfafdc7f423d7b/e2e/7322be20-897a-422d-9b9e-bb6634bde6cd.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/342d41c5bbd5653daaa7680fc7bfafdc7f423d7b/e2e/c7123d17-3c8f-4140-a76c-919c6d0db9de.md"

# Remove existing hyperlinks so they can be re-added in row order; this keeps
# the relationship ids allocated in the same A2,I2,A3,I3 order as the report
# generator.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $urlMd1, "", "", "7322be20-897a-422d-9b9e-bb6634bde6cd.md")
$zh.Hyperlinks.Add($zh.Range("I2"), $urlMd1, "", "", "7322be20-897a-422d-9b9e-bb6634bde6cd.md")
$zh.Hyperlinks.Add($zh.Range("A3"), $urlMd2, "", "", "c7123d17-3c8f-4140-a76c-919c6d0db9de.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $urlMd2, "", "", "c7123d17-3c8f-4140-a76c-919c6d0db9de.md")

# Latest Handback File (J) + Latest Handback DateTime (K)
$zh.Range("J2").Value = "7322be20-897a-422d-9b9e-bb6634bde6cd.af4b97364aa3ec483d7e16d5adc665ab4a227d0d.zh-cn.xlf"
$zh.Range("J3").Value = "c7123d17-3c8f-4140-a76c-919c6d0db9de.6c6b13daa2c6961774c24e767d0944cebcdcc3f6.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-05 04:34:14"
$zh.Range("K3").Value = "2016-09-05 04:34:14"

# Latest Target File / Latest Handback File columns now hold long file names
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $urlMd1, "", "", "7322be20-897a-422d-9b9e-bb6634bde6cd.md")
$de.Hyperlinks.Add($de.Range("I2"), $urlMd1, "", "", "7322be20-897a-422d-9b9e-bb6634bde6cd.md")
$de.Hyperlinks.Add($de.Range("A3"), $urlMd2, "", "", "c7123d17-3c8f-4140-a76c-919c6d0db9de.md")
$de.Hyperlinks.Add($de.Range("I3"), $urlMd2, "", "", "c7123d17-3c8f-4140-a76c-919c6d0db9de.md")

$de.Range("J2").Value = "7322be20-897a-422d-9b9e-bb6634bde6cd.af4b97364aa3ec483d7e16d5adc665ab4a227d0d.de-de.xlf"
$de.Range("J3").Value = "c7123d17-3c8f-4140-a76c-919c6d0db9de.6c6b13daa2c6961774c24e767d0944cebcdcc3f6.de-de.xlf"
$de.Range("K2").Value = "2016-09-05 04:34:22"
$de.Range("K3").Value = "2016-09-05 04:34:22"

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated"
